# Update the scanning-records worksheet:
#  - mark scannedStatus for several existing rows (Yes/No)
#  - fill in scan Date for rows that were previously blank
#  - append 4 new rows for trees 5502(1), 5502(2), 5543(1), 5543(2)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Excel date serial numbers (1899 date system) for the relevant scan dates.
$d0922 = 45922
$d0924 = 45924
$d0925 = 45925

# --- Rows whose scannedStatus flips from "No" to "Yes" ---
$ws.Range("E3").Value  = "Yes"
$ws.Range("E5").Value  = "Yes"
$ws.Range("E7").Value  = "Yes"
$ws.Range("E8").Value  = "Yes"
$ws.Range("E10").Value = "Yes"
$ws.Range("E11").Value = "Yes"
$ws.Range("E12").Value = "Yes"

# --- Rows 14-17: scannedStatus was blank, now "Yes" ---
$ws.Range("E14").Value = "Yes"
$ws.Range("E15").Value = "Yes"
$ws.Range("E16").Value = "Yes"
$ws.Range("E17").Value = "Yes"

# --- Row 18: scannedStatus now "No", and scan date moves to 9/25 ---
$ws.Range("E18").Value = "No"
$ws.Range("F18").Value = $d0925

# --- Rows 19-25: fill in the scan date (previously blank) ---
$ws.Range("F19").Value = $d0925
$ws.Range("F20").Value = $d0925
$ws.Range("F21").Value = $d0925
$ws.Range("F22").Value = $d0925
$ws.Range("F23").Value = $d0925
$ws.Range("F24").Value = $d0925
$ws.Range("F25").Value = $d0925

# --- New rows 26-29 for newly scanned trees ---
# Insert the new rows (shifting down, pulling formatting from the row
# above) so the date cells in column F reuse the existing date style
# instead of Excel creating a brand new cell style.
$ws.Rows("26:29").Insert(-4121, 0)  # xlShiftDown, xlFormatFromLeftOrAbove

$ws.Range("A26").Value = "AV06"
$ws.Range("B26").Value = "5502(1)"
$ws.Range("C26").Value = "ABAM"
$ws.Range("D26").Value = 60
$ws.Range("F26").Value = $d0925

$ws.Range("A27").Value = "AV06"
$ws.Range("B27").Value = "5502(2)"
$ws.Range("C27").Value = "ABAM"
$ws.Range("D27").Value = 60
$ws.Range("F27").Value = $d0925

$ws.Range("A28").Value = "AV06"
$ws.Range("B28").Value = "5543(1)"
$ws.Range("C28").Value = "PSME"
$ws.Range("D28").Value = 60
$ws.Range("F28").Value = $d0925

$ws.Range("A29").Value = "AV06"
$ws.Range("B29").Value = "5543(2)"
$ws.Range("C29").Value = "PSME"
$ws.Range("D29").Value = 60
$ws.Range("F29").Value = $d0925

# --- Update the view so it matches the author's final scroll/selection state ---
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F18:F29").Select()
